$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("article"); this pushes
# article/price/pictures one column to the right (C->D, D->E, E->F)
# and copies the formatting (incl. width) of the old column C into the
# new column, matching the widened custom-width range B:E seen in the diff.
$ws.Columns("C").Insert()

# Header for the newly inserted column
$ws.Range("C1").Value = "description"
$ws.Range("C1").Font.Bold = $true

# Description values for each product row
$ws.Range("C2").Value = "Just some beef"
$ws.Range("C3").Value = "Test description for milk"
$ws.Range("C4").Value = "Yoooo….. Ghurt."

# Leave the final selection on C4, matching the saved workbook state
$ws.Range("C4").Select()
